$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: was =1/0 (#DIV/0!) -> now =0/0 (still #DIV/0!)
$ws.Range("A1").Formula = "=0/0"

# A2 (new): array-ish formula calling an unrecognised "_xludf.NA" function,
# which Excel reports as #NAME?
$ws.Range("A2").Formula = "=_xludf.NA()"

# A3 (new): "a"+0 -> #VALUE!
$ws.Range("A3").Formula = "=""a""+0"

# Move the active selection to A4, matching the saved sheet view state
$ws.Range("A4").Select() | Out-Null
